$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.988.98"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.364.58"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'302.40"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'95.76"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "'0.503"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'33.91"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'0.124"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "'18.37"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.72"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.729.65"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "2.360.56"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'0.794"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "42.949.79"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'11.83"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "'67.93"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'234.95"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  -4.73%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "'24.58"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").Value = "'31.53"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "'17.29"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "'0.0718"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'129.50"
$ws.Range("E35").Value = "  -22.17%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.104"
$ws.Range("E36").Value = "  +3.88%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.84"
$ws.Range("E37").Value = "  +3.84%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.34"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'21.38"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").Value = "1.935.81"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.71"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.14"
$ws.Range("E47").Value = "  -9.92%  "
$ws.Range("D48").Value = "2.587.70"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'71.38"
$ws.Range("E51").Value = "  -0.96%  "

Write-Host "Applied 112 cell updates"
